$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Saldo for first data row (ANILSON) from 170129.36 to 154097.32
$ws.Range("C2").Value = 154097.32

# Remove the six rows for MARCUS, LEVI, BRASFORT, NATALIA, CLAUDIO, ANA (rows 4-9)
$ws.Range("A4:C9").EntireRow.Delete()

# Insert a single replacement row at the same position
$ws.Range("A4:C4").EntireRow.Insert()
$ws.Range("A4").Value = "'004207278"
$ws.Range("B4").Value = "CESAR"
$ws.Range("C4").Value = 9176.22
